$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows that newly get "CERO" in column N (translation marked as zero / not needed)
$ceroRows = @(13, 15, 16, 23, 25, 30)
foreach ($r in $ceroRows) {
    $ws.Cells.Item($r, 14).Value = "CERO"
}

# Row 28: new label 182 added (Asado s/h 9C) - translation still pending ("FALTA")
$ws.Cells.Item(28, 13).Value = 182
$ws.Cells.Item(28, 14).Value = "FALTA"

# Update the view: scroll the frozen pane down so row 16 is the first
# visible row below the frozen header, and leave L31 as the active selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("L31").Select()
